$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rename contraception intervention names in column C ("New intervention name").
# Written in the order that reproduces the target shared-string table ordering:
# Condoms, Female condoms, Oral contraceptive pill, Injectable, IUD, Implant,
# Female sterilization, Male sterilization.
$ws1.Range("C123").Value = "Modern contraceptive: Condoms"
$ws1.Range("C124").Value = "Modern contraceptive: Female condoms"
$ws1.Range("C122").Value = "Modern contraceptive: Oral contraceptive pill"
$ws1.Range("C126").Value = "Modern contraception: Long-acting injectable hormones"
$ws1.Range("C127").Value = "Modern contraception: Intrauterine device (IUD)"
$ws1.Range("C128").Value = "Modern contraception: Levonorgestrel-releasing implant"
$ws1.Range("C129").Value = "Modern contraception: Female sterilization"
$ws1.Range("C130").Value = "Modern contraception: Male sterilization"

# Update the visible selection on Sheet1 to match the author's final cursor position.
$ws1.Range("C127").Select() | Out-Null

# Add a new, empty Sheet2 after Sheet1.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Range("P16").Select() | Out-Null

# Restore Sheet1 as the active/selected sheet.
$ws1.Activate() | Out-Null
